$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 610.87177
$ws.Range("J17").Value = 621.6842
$ws.Range("L17").Value = 1865.0526
$ws.Range("N17").Value = -2201.0526
$ws.Range("H20").Value = 842.75
$ws.Range("I20").Value = 842.75
$ws.Range("K20").Value = 842.75
$ws.Range("M20").Value = -612.75
$ws.Range("H35").Value = 842.75
$ws.Range("I35").Value = 842.75
$ws.Range("K35").Value = 842.75
$ws.Range("M35").Value = -463.75
$ws.Range("H125").Value = 1594.1875
$ws.Range("I125").Value = 391
$ws.Range("J125").Value = 1766.0714
$ws.Range("K125").Value = 3519
$ws.Range("L125").Value = 15894.6426
$ws.Range("M125").Value = -1059
$ws.Range("N125").Value = -20814.6426
$ws.Range("H134").Value = 55071.43
$ws.Range("J134").Value = 55071.43
$ws.Range("L134").Value = 55071.43
$ws.Range("N134").Value = -65211.43

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2599.4707
$ws.Range("I61").Value = 2053.1538
$ws.Range("J61").Value = 4375
$ws.Range("K61").Value = 2053.1538
$ws.Range("L61").Value = 4375
$ws.Range("M61").Value = -1841.1538
$ws.Range("N61").Value = -4799
$ws.Range("H136").Value = 2599.4707
$ws.Range("I136").Value = 2053.1538
$ws.Range("J136").Value = 4375
$ws.Range("K136").Value = 6159.4614
$ws.Range("L136").Value = 13125
$ws.Range("M136").Value = -3609.4614
$ws.Range("N136").Value = -18225

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 80516.22
$ws.Range("J82").Value = 20626.4
$ws.Range("L82").Value = 20626.4
$ws.Range("N82").Value = -21392.4
$ws.Range("H85").Value = 80516.22
$ws.Range("J85").Value = 20626.4
$ws.Range("L85").Value = 20626.4
$ws.Range("N85").Value = -23278.4
$ws.Range("H134").Value = 133127.56
$ws.Range("I134").Value = 233964.31
$ws.Range("J134").Value = 2039.8
$ws.Range("K134").Value = 701892.9299999999
$ws.Range("L134").Value = 6119.4
$ws.Range("M134").Value = -699357.9299999999
$ws.Range("N134").Value = -11189.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2358.2
$ws.Range("I58").Value = 1854.5714
$ws.Range("J58").Value = 3533.3333
$ws.Range("K58").Value = 1854.5714
$ws.Range("L58").Value = 3533.3333
$ws.Range("M58").Value = -1651.5714
$ws.Range("N58").Value = -3939.3333
$ws.Range("H136").Value = 2358.2
$ws.Range("I136").Value = 1854.5714
$ws.Range("J136").Value = 3533.3333
$ws.Range("K136").Value = 5563.7142
$ws.Range("L136").Value = 10599.9999
$ws.Range("M136").Value = -3013.7142
$ws.Range("N136").Value = -15699.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9525.429
$ws.Range("I3").Value = 7838
$ws.Range("J3").Value = 13744
$ws.Range("K3").Value = 23514
$ws.Range("L3").Value = 41232
$ws.Range("M3").Value = -23402
$ws.Range("N3").Value = -41456
$ws.Range("H5").Value = 1212973.9
$ws.Range("I5").Value = 680
$ws.Range("J5").Value = 2667726.5
$ws.Range("K5").Value = 2040
$ws.Range("L5").Value = 8003179.5
$ws.Range("M5").Value = -1928
$ws.Range("N5").Value = -8003403.5
$ws.Range("H64").Value = 7849.647
$ws.Range("I64").Value = 512
$ws.Range("J64").Value = 8308.25
$ws.Range("K64").Value = 1536
$ws.Range("L64").Value = 24924.75
$ws.Range("M64").Value = -1266
$ws.Range("N64").Value = -25464.75
$ws.Range("H67").Value = 7849.647
$ws.Range("I67").Value = 512
$ws.Range("J67").Value = 8308.25
$ws.Range("K67").Value = 1536
$ws.Range("L67").Value = 24924.75
$ws.Range("M67").Value = -600
$ws.Range("N67").Value = -26796.75
$ws.Range("H68").Value = 1233.3334
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1233.3334
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").Value = 3700.0002
$ws.Range("N68").Value = -5322.0002
$ws.Range("H71").Value = 1233.3334
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1233.3334
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").Value = 11100.0006
$ws.Range("N71").Value = -19212.0006
$ws.Range("H107").Value = 988.1818
$ws.Range("I107").Value = 482.66666
$ws.Range("J107").Value = 1177.75
$ws.Range("K107").Value = 1447.99998
$ws.Range("L107").Value = 3533.25
$ws.Range("M107").Value = 472.0000199999999
$ws.Range("N107").Value = -7373.25
$ws.Range("H123").Value = 5250
$ws.Range("I123").Value = 1862.5
$ws.Range("K123").Value = 5587.5
$ws.Range("M123").Value = -3137.5
$ws.Range("H135").Value = 1212973.9
$ws.Range("I135").Value = 680
$ws.Range("J135").Value = 2667726.5
$ws.Range("K135").Value = 6120
$ws.Range("L135").Value = 24009538.5
$ws.Range("M135").Value = -3585
$ws.Range("N135").Value = -24014608.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21385.793
$ws.Range("I70").Value = 36251.668
$ws.Range("J70").Value = 5033.3335
$ws.Range("K70").Value = 36251.668
$ws.Range("L70").Value = 5033.3335
$ws.Range("M70").Value = -35981.668
$ws.Range("N70").Value = -5573.3335
$ws.Range("H73").Value = 21385.793
$ws.Range("I73").Value = 36251.668
$ws.Range("J73").Value = 5033.3335
$ws.Range("K73").Value = 36251.668
$ws.Range("L73").Value = 5033.3335
$ws.Range("M73").Value = -35315.668
$ws.Range("N73").Value = -6905.3335
$ws.Range("H132").Value = 2526
$ws.Range("I132").Value = 2284
$ws.Range("J132").Value = 3044.5715
$ws.Range("K132").Value = 6852
$ws.Range("L132").Value = 9133.7145
$ws.Range("M132").Value = -4322
$ws.Range("N132").Value = -14193.7145
$ws.Range("H135").Value = 72853.336
$ws.Range("J135").Value = 72853.336
$ws.Range("L135").Value = 72853.336
$ws.Range("N135").Value = -82993.336

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3706087
$ws.Range("I7").Value = 5883980
$ws.Range("J7").Value = 3669.3
$ws.Range("K7").Value = 5883980
$ws.Range("L7").Value = 3669.3
$ws.Range("M7").Value = -5883868
$ws.Range("N7").Value = -3893.3
$ws.Range("H126").Value = 3706087
$ws.Range("I126").Value = 5883980
$ws.Range("J126").Value = 3669.3
$ws.Range("K126").Value = 17651940
$ws.Range("L126").Value = 11007.9
$ws.Range("M126").Value = -17649470
$ws.Range("N126").Value = -15947.9
$ws.Range("H132").Value = 2254.6924
$ws.Range("I132").Value = 1861.1666
$ws.Range("J132").Value = 2592
$ws.Range("K132").Value = 5583.4998
$ws.Range("L132").Value = 7776
$ws.Range("M132").Value = -3053.4998
$ws.Range("N132").Value = -12836
$ws.Range("H133").Value = 51136
$ws.Range("J133").Value = 51136
$ws.Range("L133").Value = 51136
$ws.Range("N133").Value = -56196
$ws.Range("H136").Value = 1529.9656
$ws.Range("I136").Value = 1230.1818
$ws.Range("J136").Value = 2472.1428
$ws.Range("K136").Value = 3690.5454
$ws.Range("L136").Value = 7416.428400000001
$ws.Range("M136").Value = -1140.5454
$ws.Range("N136").Value = -12516.4284

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 60000
$ws.Range("J46").Value = 60000
$ws.Range("L46").Value = 60000
$ws.Range("N46").Value = -60462
$ws.Range("H81").Value = 5744.913
$ws.Range("I81").Value = 10613.1
$ws.Range("J81").Value = 2000.1538
$ws.Range("K81").Value = 21226.2
$ws.Range("L81").Value = 4000.3076
$ws.Range("M81").Value = -20165.2
$ws.Range("N81").Value = -6122.3076
$ws.Range("H84").Value = 5744.913
$ws.Range("I84").Value = 10613.1
$ws.Range("J84").Value = 2000.1538
$ws.Range("K84").Value = 106131
$ws.Range("L84").Value = 20001.538
$ws.Range("M84").Value = -100827
$ws.Range("N84").Value = -30609.538
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 180000
$ws.Range("N134").Value = -185070
$ws.Range("H135").Value = 57295.715
$ws.Range("J135").Value = 57295.715
$ws.Range("L135").Value = 57295.715
$ws.Range("N135").Value = -67435.715
$ws.Range("H136").Value = 1923.3462
$ws.Range("I136").Value = 1566.3914
$ws.Range("J136").Value = 4660
$ws.Range("K136").Value = 4699.174199999999
$ws.Range("L136").Value = 13980
$ws.Range("M136").Value = -2149.174199999999
$ws.Range("N136").Value = -19080
$ws.Range("H141").Value = 55350
$ws.Range("L141").Value = 55350
$ws.Range("N141").Value = -65710
